$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the reported IFRS figures (columns D:AJ) for the annual rows
# (rows 2-6, fiscal years 2014-2018) with the revised/rescaled values.

# Row 2 - 2014/12 (IFRS연결)
$ws.Range("D2").Value = 3570
$ws.Range("E2").Value = 312
$ws.Range("F2").Value = 312
$ws.Range("G2").Value = 342
$ws.Range("H2").Value = 256
$ws.Range("I2").Value = 258
$ws.Range("J2").Value = -1
$ws.Range("K2").Value = 3524
$ws.Range("L2").Value = 897
$ws.Range("M2").Value = 2627
$ws.Range("N2").Value = 2440
$ws.Range("O2").Value = 187
$ws.Range("P2").Value = 167
$ws.Range("Q2").Value = 359
$ws.Range("R2").Value = -504
$ws.Range("S2").Value = -12
$ws.Range("T2").Value = 95
$ws.Range("U2").Value = 264
$ws.Range("V2").Value = 23
$ws.Range("W2").Value = 8.74
$ws.Range("X2").Value = 7.18
$ws.Range("Y2").Value = 11.07
$ws.Range("Z2").Value = 7.79
$ws.Range("AA2").Value = 34.15
$ws.Range("AB2").Value = 1397.23
$ws.Range("AC2").Value = 771
$ws.Range("AD2").Value = 13.57
$ws.Range("AE2").Value = 7597
$ws.Range("AF2").Value = 1.38
$ws.Range("AG2").Value = 100
$ws.Range("AH2").Value = 0.95
$ws.Range("AI2").Value = 12.47
$ws.Range("AJ2").Value = 33393384

# Row 3 - 2015/12 (IFRS연결)
$ws.Range("D3").Value = 4355
$ws.Range("E3").Value = 538
$ws.Range("F3").Value = 538
$ws.Range("G3").Value = 613
$ws.Range("H3").Value = 457
$ws.Range("I3").Value = 429
$ws.Range("J3").Value = 28
$ws.Range("K3").Value = 4039
$ws.Range("L3").Value = 1010
$ws.Range("M3").Value = 3029
$ws.Range("N3").Value = 2832
$ws.Range("O3").Value = 196
$ws.Range("P3").Value = 167
$ws.Range("Q3").Value = 491
$ws.Range("R3").Value = -506
$ws.Range("S3").Value = -32
$ws.Range("T3").Value = 237
$ws.Range("U3").Value = 254
$ws.Range("V3").Value = 28
$ws.Range("W3").Value = 12.34
$ws.Range("X3").Value = 10.49
$ws.Range("Y3").Value = 16.27
$ws.Range("Z3").Value = 12.08
$ws.Range("AA3").Value = 33.35
$ws.Range("AB3").Value = 1634.81
$ws.Range("AC3").Value = 1284
$ws.Range("AD3").Value = 11.53
$ws.Range("AE3").Value = 8832
$ws.Range("AF3").Value = 1.68
$ws.Range("AG3").Value = 150
$ws.Range("AH3").Value = 1.01
$ws.Range("AI3").Value = 11.22
$ws.Range("AJ3").Value = 33393384

# Row 4 - 2016/12 (IFRS연결)
$ws.Range("D4").Value = 2604
$ws.Range("E4").Value = 217
$ws.Range("F4").Value = 541
$ws.Range("G4").Value = 303
$ws.Range("H4").Value = 511
$ws.Range("I4").Value = 480
$ws.Range("J4").Value = 32
$ws.Range("K4").Value = 4729
$ws.Range("L4").Value = 1328
$ws.Range("M4").Value = 3401
$ws.Range("N4").Value = 3166
$ws.Range("O4").Value = 236
$ws.Range("P4").Value = 167
$ws.Range("Q4").Value = 296
$ws.Range("R4").Value = -190
$ws.Range("S4").Value = -129
$ws.Range("T4").Value = 383
$ws.Range("U4").Value = -86
$ws.Range("V4").Value = 38
$ws.Range("W4").Value = 8.34
$ws.Range("X4").Value = 19.64
$ws.Range("Y4").Value = 15.99
$ws.Range("Z4").Value = 11.67
$ws.Range("AA4").Value = 39.05
$ws.Range("AB4").Value = 1893.19
$ws.Range("AC4").Value = 1436
$ws.Range("AD4").Value = 14.62
$ws.Range("AE4").Value = 10086
$ws.Range("AF4").Value = 2.08
$ws.Range("AG4").Value = 150
$ws.Range("AH4").Value = 0.71
$ws.Range("AI4").Value = 9.82
$ws.Range("AJ4").Value = 33393384

# Row 5 - 2017/12 (IFRS연결)
$ws.Range("D5").Value = 5730
$ws.Range("E5").Value = 555
$ws.Range("F5").Value = 555
$ws.Range("G5").Value = 751
$ws.Range("H5").Value = 850
$ws.Range("I5").Value = 734
$ws.Range("J5").Value = 117
$ws.Range("K5").Value = 3761
$ws.Range("L5").Value = 1285
$ws.Range("M5").Value = 2477
$ws.Range("N5").Value = 2127
$ws.Range("O5").Value = 350
$ws.Range("P5").Value = 68
$ws.Range("Q5").Value = 562
$ws.Range("R5").Value = -470
$ws.Range("S5").Value = -112
$ws.Range("T5").Value = 251
$ws.Range("U5").Value = 310
$ws.Range("V5").Value = 76
$ws.Range("W5").Value = 9.69
$ws.Range("X5").Value = 14.84
$ws.Range("Y5").Value = 27.73
$ws.Range("Z5").Value = 20.03
$ws.Range("AA5").Value = 51.87
$ws.Range("AB5").Value = 5287.61
$ws.Range("AC5").Value = 2439
$ws.Range("AD5").Value = 10.13
$ws.Range("AE5").Value = 17092
$ws.Range("AF5").Value = 1.45
$ws.Range("AG5").Value = 180
$ws.Range("AH5").Value = 0.73
$ws.Range("AI5").Value = 3.05
$ws.Range("AJ5").Value = 13554044

# Row 6 - 2018/12 (IFRS연결)
$ws.Range("D6").Value = 5261
$ws.Range("E6").Value = 761
$ws.Range("F6").Value = 761
$ws.Range("G6").Value = 2860
$ws.Range("H6").Value = 2106
$ws.Range("I6").Value = 1947
$ws.Range("K6").Value = 5721
$ws.Range("L6").Value = 1575
$ws.Range("M6").Value = 4146
$ws.Range("N6").Value = 3650
$ws.Range("P6").Value = 68
$ws.Range("Q6").Value = 596
$ws.Range("R6").Value = -236
$ws.Range("S6").Value = 0
$ws.Range("T6").Value = 229
$ws.Range("U6").Value = 367
$ws.Range("V6").Value = 102
$ws.Range("W6").Value = 14.47
$ws.Range("X6").Value = 40.04
$ws.Range("Y6").Value = 67.41
$ws.Range("Z6").Value = 44.43
$ws.Range("AA6").Value = 37.99
$ws.Range("AB6").Value = 5925.56
$ws.Range("AC6").Value = 14365
$ws.Range("AD6").Value = 0.82
$ws.Range("AE6").Value = 29330
$ws.Range("AF6").Value = 0.4
$ws.Range("AG6").Value = 220
$ws.Range("AH6").Value = 1.87
$ws.Range("AI6").Value = 1.41
$ws.Range("AJ6").Value = 13554044

# Rows 7-9 (2019/12(E), 2020/12(E), 2021/12(E)) no longer have projected
# figures available; clear D:AJ so only A (index), B (period) and C
# (fiscal-period label) remain populated, matching the other sheets.
$ws.Range("D7:AJ9").ClearContents()
